$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$newIds = @{
    2 = "23fb7f2a-6a9c-4ea1-82f2-623ea4ba38f5"
    3 = "1a5444e6-f505-4fb4-bf54-050f84516d6e"
    4 = "f772d4ce-83c6-4eac-8470-e2a809381ee8"
    5 = "72d47d98-219f-4b69-82ac-8638befeb868"
    6 = "dbe638b1-724d-4104-9037-83d1da57449e"
    7 = "3305aeb2-9df5-4d71-9031-74c721b020d4"
    8 = "8c11661c-f148-463e-b7d4-6edcfc1a1d78"
    9 = "4379df56-8053-41f8-9c58-7ac9d1f7ba91"
    10 = "24a6496e-070c-4230-8c85-7e20e9bae756"
    11 = "7009a27f-22c6-4041-b738-00c95dafe3a6"
    12 = "fd3cd4fe-456d-4d99-af1a-4671de0a1a04"
    13 = "5e332762-b7e9-4dc1-b80c-4e7e390f282c"
    14 = "2106dd5f-bbf4-4953-8987-276acad055ae"
    15 = "ac9ce6b3-a435-4656-9479-389fdd5dd652"
    16 = "12cfa314-4c5c-45e2-a1c9-79e5d31329cb"
    17 = "e8128046-ca60-4932-95ab-56674ee7dd23"
    18 = "fe6e8e34-8aae-49ea-bc74-6a75e031b05a"
    19 = "8f1769be-b9ff-470c-b0ba-ed0d1b1c97dd"
    20 = "825edb45-da5d-4e71-b954-52f3ddc0957e"
    21 = "63b94fcd-375f-4a5a-8815-1a4d562fdbee"
    22 = "33ecd688-b39c-41c0-8519-ad86d3157257"
    23 = "b927a30a-d1dd-42cc-a515-563f9490e91e"
    24 = "7b65c2cf-0a0a-468a-9110-369ab9192dea"
    25 = "cba00a01-e1cd-46a7-a278-4c9aedbcc372"
    26 = "5dab577d-83dc-45db-ae48-793fa2301ec2"
    27 = "a0c5b76a-b6cf-4210-a317-cd268e0261f2"
    28 = "543be7af-d057-4d30-8ad9-471133f2ae6e"
    29 = "e450b46f-3e68-4f8d-ba48-54f23054f79d"
    30 = "cf3d1053-f6ce-42e6-9610-d85d69d621c4"
    31 = "4fcfd071-46b4-4f0f-98f2-745cc7d0e57d"
    32 = "6e0a2052-5d78-4b7d-9969-dbe5980ecf57"
    33 = "cbb27627-d1db-4cf4-85da-12610657dfca"
    34 = "c6813f05-fe68-47b9-b454-69facc02b939"
    35 = "3ff99f77-a629-4ffb-9dfd-0e0c386fcd00"
    36 = "5b837cfa-2c99-4d99-9c70-cf397658e377"
    37 = "1411cf27-3598-47c6-8093-3d99e5e3a8fa"
    38 = "4921d556-017b-42c3-922c-0849ff8b7553"
    39 = "77836f79-1d43-4f17-96a5-7230d1a82122"
    40 = "d8054ee3-8c5f-46a5-9723-00551f37feed"
    41 = "381b5d4c-ee98-4bb1-b4f8-a7363639ea46"
    42 = "0cf8e999-ff7a-4dda-ae96-2760514c1cb7"
    43 = "cc375c51-c152-4626-8662-eba878668df7"
    44 = "57500dbf-377f-40d2-976a-c561e40cb632"
    45 = "dbde23d4-d66a-4820-8603-da4400657570"
    46 = "71041473-44ad-4cdd-85e6-c297e8c088b7"
    47 = "8fff461b-93a1-42e8-8e81-07c78e81a043"
    48 = "d20c7a75-c109-4792-88be-22a70256e2fe"
    49 = "2cd08745-505c-4718-987a-5f7c89bc05a1"
    50 = "568fcd18-4265-4d63-a08d-d6cb3320eebc"
    51 = "e28a1b92-816c-4102-a9c2-0c51db88e7ae"
    52 = "eeb19dbd-7a50-44d5-816f-d85c03f03876"
    53 = "855ad8d6-d82c-4912-82bd-c8f7ae6b6c15"
    54 = "fd93f280-5bf9-4d3c-8802-5dcc3ec8dbb9"
    55 = "b9ba049f-6f75-413f-8d3c-e15be59e1bd2"
    56 = "4955f804-2ad0-4e3b-a76b-ca09e89137cf"
    57 = "ff3c1f34-469e-40fb-bdc4-4980ec70eed9"
    58 = "41f4b5aa-01b7-46ed-951c-622e2adc39a3"
    59 = "da391883-92d9-4a36-b3aa-32d33ac6e15b"
    60 = "d8c667b3-0ec2-4c86-b0f0-e545a1c4dad5"
    61 = "849304f8-3de8-4ec9-b685-c95d4d00dde6"
    62 = "bebf1abc-dc5e-432e-8eb1-5281b628cb7e"
    63 = "5e5b90e3-bb6b-4a43-bbd6-56aeb720d6fe"
    64 = "6a580fcc-b6a9-49c0-b317-937f70b65f2b"
    65 = "ebc973e3-c637-45a0-94c8-e54b1c8daaea"
    66 = "3275935e-7408-42ea-8bca-d032bf5d82ad"
    67 = "28a49d02-f4e0-475e-b5b9-f247fcff9329"
    68 = "fe28baaa-e512-4918-8ffb-07b1632fb952"
    69 = "be5d723e-5b6c-4748-a0c0-25715493dae1"
    70 = "e3d5fdf7-c58a-4dc8-b2c6-f59bbe71ea72"
    71 = "305d8d3a-d319-486f-8fcc-0a2c3f26efd3"
    72 = "262e316a-5a33-4b8d-8a5a-2a6254bfde44"
    73 = "61304178-e0a9-46de-99c4-247de83a1c6c"
    74 = "b9bd9660-acc9-4098-8b3c-737824dfeea0"
    75 = "7b5a0490-acc4-4399-9e49-9c4883e89571"
    76 = "ba0ef346-179e-445c-b01a-40c02d31a576"
    77 = "8c255a5a-c62c-46cd-b507-5561cc50578f"
    78 = "bddbf100-5284-488a-a8e3-a4c2a9c437e9"
    79 = "f68585c2-56cb-47b6-87c5-b53e9a92ede2"
    80 = "de7c5407-bd0d-49b7-b76b-8191f6869804"
    96 = "353678b9-9a62-4dab-bdc0-f83a58107132"
}

foreach ($row in $newIds.Keys) {
    $ws.Cells.Item($row, 12).Value = $newIds[$row]
}

